$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Week 40 (rows 42-48): Rief (column C) attended Monday (row 46) for 2 hours ---
$ws.Range("C46").Value = 2

# --- Week 43 (rows 58-64): Walter(D), Paige(E), Rick(G), Benjamin(H), Hadewij(I)
#     attended Friday (row 63) for 4 hours each ---
$ws.Range("D63").Value = 4
$ws.Range("E63").Value = 4
$ws.Range("G63").Value = 4
$ws.Range("H63").Value = 4
$ws.Range("I63").Value = 4

# --- Week 43 total Game-Lab hours per week (row 64) bumped from 10 to 18 ---
$ws.Range("B64").Value = 18

# --- Move the on-screen selection to reflect where the author was working ---
$ws.Range("J37").Select() | Out-Null

$wb.Application.Calculate() | Out-Null
